$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Delete the whole paragraph "Thompson and Gonzalez (2017) need
#    movement + future climate for predictions, can't do with data
#    scarcity" (it is being dropped entirely).
# ------------------------------------------------------------------
$apos = [char]0x2019
$thompsonPara = $d.Paragraphs.Item(4)
$thompsonText = "Thompson and Gonzalez (2017) need movement + future climate for predictions, can" + $apos + "t do with data scarcity"
if ($thompsonPara.Range.Text.TrimEnd([char]13) -eq $thompsonText) {
    $thompsonPara.Range.Delete()
} else {
    throw "Paragraph 4 did not match expected Thompson/Gonzalez text: [$($thompsonPara.Range.Text)]"
}

# ------------------------------------------------------------------
# 2. Extend the first body paragraph: replace the closing period after
#    "...Trojelsgaard and Olesen 2016)" with the new continuation
#    sentence that leads into the (now-adjacent) "(Hui and Richardson
#    2019; Guiden et al. 2019)" citation paragraph.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$oldTail = "(Pellissier et al. 2017; Tr" + [char]0x00F8 + "jelsgaard and Olesen 2016)."
$newTail = "(Pellissier et al. 2017; Tr" + [char]0x00F8 + "jelsgaard and Olesen 2016). As ecosystems and climates are changing rapidly, ecologists realized that networks are at risk or unravelling, being invaded by exotic species that can destabilize them, or adopt entirely novel configurations "
$found = $find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)
if (-not $found) {
    throw "Could not find the Pellissier/Trojelsgaard sentence to extend."
}

# ------------------------------------------------------------------
# 3. Insert the remainder of the new material right after the
#    "(Hui and Richardson 2019; Guiden et al. 2019)" citation, which
#    currently still lives in its own paragraph.
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$oldHui = "(Hui and Richardson 2019; Guiden et al. 2019)"
$newHui = "(Hui and Richardson 2019; Guiden et al. 2019). Simulation studies seem to suggest that knowing the shape of the extant network is not sufficient (Thompson and Gonzalez 2017), and that it needs to be supplemented by additional data on species properties, climate, and climate projection."
$found2 = $find2.Execute($oldHui, $true, $false, $false, $false, $false, $true, 1, $false, $newHui, 2)
if (-not $found2) {
    throw "Could not find the Hui/Richardson citation to extend."
}

# ------------------------------------------------------------------
# 4. Merge the (Hui and Richardson...) paragraph back into the
#    preceding paragraph so the whole passage forms a single
#    paragraph, by deleting the paragraph mark that separates them.
# ------------------------------------------------------------------
$ecoPara = $d.Paragraphs.Item(3)
$ecoRange = $ecoPara.Range
$markRange = $d.Range($ecoRange.End - 1, $ecoRange.End)
$markRange.Delete()

# Deleting the paragraph mark can make the merged paragraph inherit the
# style of the following (now-consumed) paragraph mark; force it back to
# the original "FirstParagraph" style used by this introductory paragraph.
$d.Paragraphs.Item(3).Range.Style = "FirstParagraph"

Write-Output "done"
